# Generate Report for Handback
# Updates the zh-cn and de-de localization-status worksheets with the
# latest handback info for 8a45673f-20d6-4ebb-bb2b-8fae5966994c (row 7):
#  - marks the target/handback file + handback datetime
#  - records an Error Detail because the handback file version is stale
#  - widens the "Error Detail" column (P) so the message is readable
#  - adds a hyperlink on the new "Latest Target File" cell (I7)

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/edc3ce23a54c34bbaf4eb31160f6fd2fa98d5f78/e2e/8a45673f-20d6-4ebb-bb2b-8fae5966994c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0f34e047848b3e6bc104bd4b0e8a91c96e2b225f/e2e/8a45673f-20d6-4ebb-bb2b-8fae5966994c.md."
$targetFileUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0f34e047848b3e6bc104bd4b0e8a91c96e2b225f/e2e/8a45673f-20d6-4ebb-bb2b-8fae5966994c.md"
$mdName = "8a45673f-20d6-4ebb-bb2b-8fae5966994c.md"

$sheetNames = @("zh-cn", "de-de")
$handbackDates = @{ "zh-cn" = "2016-08-23 10:44:43"; "de-de" = "2016-08-23 10:44:50" }

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Widen the Error Detail column (P) to fit the new message.
    $ws.Columns.Item(16).ColumnWidth = 39.17

    # Row 7 is the 8a45673f-20d6-4ebb-bb2b-8fae5966994c.md entry.
    $handoffFile = $ws.Range("G7").Text

    $ws.Range("I7").Value = $mdName
    $ws.Range("I7").Font.Underline = 2
    $ws.Range("I7").Font.Color = 15570276
    $ws.Hyperlinks.Add($ws.Range("I7"), $targetFileUrl, [Type]::Missing, [Type]::Missing, $mdName)

    $ws.Range("J7").Value = $handoffFile
    $ws.Range("K7").Value = $handbackDates[$sheetName]
    $ws.Range("P7").Value = $errorDetail
}
